# Updated cryptos list - apply per-cell text updates from the scrape diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.515.15"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.728.62"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4816"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2674"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.731.99"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07191"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.59"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.529"
$ws.Range("D14").ClearFormats()
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.525.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006953"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.953.52"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.528"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.821"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.260"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.777"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.978"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08029"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.700"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04519"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6267"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9101"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.074"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.386"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.24"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -9.84%  "
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.556"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.975"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.80%  "
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05366"
$ws.Range("D47").ClearFormats()
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.842"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3406"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.54%  "
